$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New config rows: URL, Username (as hyperlink), Password
$ws.Range("B8").Value = "http://bws-selfservice.unileversolutions.com/"

$ws.Range("B9").Value = "skumar213@sapient.com"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:skumar213@sapient.com")

$ws.Range("B10").Value = "334512s"

# Switch default browser selection from Chrome to IE
$ws.Range("B3").Value = "IE"

# Update selected cell
[void]$ws.Range("B3").Select()
